# The commit renames the workbook's only worksheet from the default
# "Sheet1" to "sort_table_SP" (reflecting that supplemental file 5 holds
# the sorted SP drug-class prediction table). Renaming the sheet via the
# Excel object model automatically keeps everything that references the
# sheet by name in sync - in particular the workbook-scoped
# "_xlnm._FilterDatabase" defined name created by the AutoFilter on this
# sheet, whose reference switches from "Sheet1!$A$1:$F$282" to
# "sort_table_SP!$A$1:$F$282".

$wb = $excel.ActiveWorkbook

# Grab the sheet by its current (default) name and rename it.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "sort_table_SP"
